{"js": "// Apply the radicado/date refresh + requester re-assignment + subject/body\n// wording edits described by the commit.\n//\n// Strategy: each change is a precise, context-scoped literal text\n// replacement driven by Body.search (Office.js \"Find & Replace\" surface),\n// mirroring exactly what a human reviewer would have typed.\n\nconst body = context.document.body;\n\n// Pairs of [searchText, replacementText]. Order does not matter since each\n// search string is unique/unambiguous in the document, except the radicado\n// number which legitimately repeats twice (bare number + \"Radicado: \"\n// line) and must become the new number in both places.\nconst replacements = [\n  // Radicado number (appears twice: barcode line + \"Radicado: \" line).\n  [\"20210202183246\", \"20210311155323\"],\n  // Long-form Spanish date/time stamp.\n  [\"Bogot\u00e1 / D.C -  2021/02/02 18:32:46\", \"Bogot\u00e1 / D.C -  2021/03/11 15:53:23\"],\n  // Addressee block: name, address, email.\n  [\"Donaldo Jinete Forero\", \"Juan Perez Martinez\"],\n  [\"Cra 46 No 123 -66 Apt 502 - Bogot\u00e1 / D.C\", \"Cra 12 No 23-45 - Bogot\u00e1 / D.C\"],\n  [\"donaldo.jinette@gmail.com\", \"juan.perez@gmail.com\"],\n  // Subject line.\n  [\"RESPUESTA Hoja de vida del funcionario\", \"RESPUESTA Cat\u00e1logo de sistemas de informaci\u00f3n\"],\n  // Body opening sentence.\n  [\"Gracias por comunicarse con nosotros\", \"Muchas gracias por la respuesta\"],\n  // Signature block name.\n  [\"Donaldo Rafael Jinete Forero \", \"Miguel Cubides \"],\n];\n\nfor (const [searchText, replacementText] of replacements) {\n  const results = body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replacementText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// The trailing \" MinCiencias\" run that used to follow \"Sistema RINO\" is\n// dropped entirely (not just its text blanked out), so erase that literal\n// substring \u2014 this removes the whole phrase, collapsing the paragraph back\n// to just \"Sistema RINO\".\nconst minCienciasResults = body.search(\" MinCiencias\", { matchCase: true });\nminCienciasResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < minCienciasResults.items.length; i++) {\n  minCienciasResults.items[i].insertText(\"\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Apply the radicado/date refresh + requester re-assignment + subject/body\n# wording edits described by the commit.\n#\n# Strategy: each change is a precise, context-scoped literal text\n# replacement driven by Find/Replace over the whole document content,\n# mirroring exactly what a human reviewer would have typed in Word.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($find, $replace) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Execute(\n        $find,       # FindText\n        $true,       # MatchCase\n        $false,      # MatchWholeWord\n        $false,      # MatchWildcards\n        $false,      # MatchSoundsLike\n        $false,      # MatchAllWordForms\n        $true,       # Forward\n        1,           # Wrap = wdFindContinue\n        $false,      # Format\n        $replace,    # ReplaceWith\n        2            # Replace = wdReplaceAll\n    )\n}\n\n# Radicado number (appears twice: barcode line + \"Radicado: \" line).\nReplace-Text \"20210202183246\" \"20210311155323\"\n\n# Long-form Spanish date/time stamp.\nReplace-Text \"Bogot\u00e1 / D.C -  2021/02/02 18:32:46\" \"Bogot\u00e1 / D.C -  2021/03/11 15:53:23\"\n\n# Addressee block: name, address, email.\nReplace-Text \"Donaldo Jinete Forero\" \"Juan Perez Martinez\"\nReplace-Text \"Cra 46 No 123 -66 Apt 502 - Bogot\u00e1 / D.C\" \"Cra 12 No 23-45 - Bogot\u00e1 / D.C\"\nReplace-Text \"donaldo.jinette@gmail.com\" \"juan.perez@gmail.com\"\n\n# Subject line.\nReplace-Text \"RESPUESTA Hoja de vida del funcionario\" \"RESPUESTA Cat\u00e1logo de sistemas de informaci\u00f3n\"\n\n# Body opening sentence.\nReplace-Text \"Gracias por comunicarse con nosotros\" \"Muchas gracias por la respuesta\"\n\n# Signature block name (trailing space is part of the original run text).\nReplace-Text \"Donaldo Rafael Jinete Forero \" \"Miguel Cubides \"\n\n# The trailing \" MinCiencias\" run that used to follow \"Sistema RINO\" is\n# dropped entirely, so erase that literal substring \u2014 this removes the\n# whole phrase, collapsing the paragraph back to just \"Sistema RINO\".\nReplace-Text \" MinCiencias\" \"\"\n"}
